$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 113, shifting existing rows 113:178 down to 114:179
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new weekly price observation
$ws.Cells.Item(113, 1).Value = 4
$ws.Cells.Item(113, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(113, 3).Value = "Los Lagos"
$ws.Cells.Item(113, 4).Value = 44606
$ws.Cells.Item(113, 5).Value = 10
$ws.Cells.Item(113, 6).Value = 100112039
$ws.Cells.Item(113, 7).Value = "Ciboulette"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 80
$ws.Cells.Item(113, 11).Value = 2500
$ws.Cells.Item(113, 12).Value = 2500
$ws.Cells.Item(113, 13).Value = 2500
$ws.Cells.Item(113, 14).Value = "`$/docena de atados"
$ws.Cells.Item(113, 15).Value = "Región Metropolitana"
$ws.Cells.Item(113, 16).Value = 833
$ws.Cells.Item(113, 17).Value = 3
$ws.Cells.Item(113, 18).Value = "Hortaliza"
